# Add a new "Lifetime (yr)" row to the equipment-lifetime table on Sheet1.
# This mirrors the new BioSTEAM equipment-lifetime feature: every piece of
# equipment gets an extra boolean flag row, just like the existing
# "Bag unloader" row (row 12) above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = 12
$newRow = 13

# A13: label cell - reuse the same formatting as A12 (bordered, left aligned)
$ws.Range("A$lastRow").Copy() | Out-Null
$ws.Range("A$newRow").PasteSpecial(-4122) | Out-Null
$ws.Range("A$newRow").Value = "Lifetime (yr)"

# B13:AN13: boolean flag cells - reuse the same formatting as B12:AN12
# (centered, no border) and default every flag to FALSE.
$ws.Range("B${lastRow}:AN${lastRow}").Copy() | Out-Null
$ws.Range("B${newRow}:AN${newRow}").PasteSpecial(-4122) | Out-Null
$ws.Range("B${newRow}:AN${newRow}").Value = $false

$excel.CutCopyMode = $false

# Move the active selection to the newly added row, as in the saved workbook.
$ws.Range("A$newRow").Select()
